$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Periodo Mora" column (E) values for rows 16-24 are reversed in order,
# and the "Valor Mora" column (F) keeps 40000 for all periods except the
# newest one (2103, now on row 16) which carries the 30666 value that
# previously belonged to the last row (old EC entries replaced by new ones).
$periodos = @("2103", "2102", "2101", "2012", "2011", "2010", "2009", "2008", "2007")
$valores  = @(30666, 40000, 40000, 40000, 40000, 40000, 40000, 40000, 40000)

for ($i = 0; $i -lt 9; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periodos[$i]
    $ws.Cells.Item($row, 6).Value = $valores[$i]
}
